$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 183, shifting existing rows 183-202 down to 184-203
$ws.Rows("183:183").Insert()

# Populate the newly inserted row 183 with the new weekly price record
$ws.Range("A183").Value = 10
$ws.Range("B183").Value = "Vega Modelo de Temuco"
$ws.Range("C183").Value = "La Araucanía"
$ws.Range("D183").Value = 45154
$ws.Range("E183").Value = 9
$ws.Range("F183").Value = 100112035
$ws.Range("G183").Value = "Bruselas (repollito)"
$ws.Range("H183").Value = "Sin especificar"
$ws.Range("I183").Value = "Primera"
$ws.Range("J183").Value = 30
$ws.Range("K183").Value = 25000
$ws.Range("L183").Value = 25000
$ws.Range("M183").Value = 25000
$ws.Range("N183").Value = "$/malla 15 kilos"
$ws.Range("O183").Value = "Región Metropolitana"
$ws.Range("P183").Value = 1667
$ws.Range("Q183").Value = 15
$ws.Range("R183").Value = "Hortaliza"
